$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D5").Value = "Talk"
$ws.Range("F5").Value = 109

$ws.Rows.Item(6).Delete()

$ws.Range("F5").Select() | Out-Null
